$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I19").Value = -0.2080064881564166
$ws.Range("J19").Value = 0.1304043896348948
$ws.Range("K19").Value = -0.09342636787466206
$ws.Range("L19").Value = 1.993387849628447
